$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (column A, row 1 uses the bold/standard style)
# onto the new column F range (F1:F21), then fill in the values.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:F21").PasteSpecial(-4122) | Out-Null

# Match the shared-string interning order of the original edit: the
# "Yes" outlier (F16) was entered first, then the "MB Endpoint" header,
# then the remaining "No" values.
$ws.Range("F16").Value = "Yes"
$ws.Range("F1").Value = "MB Endpoint"

$ws.Range("F2").Value = "No"
$ws.Range("F3").Value = "No"
$ws.Range("F4").Value = "No"
$ws.Range("F5").Value = "No"
$ws.Range("F6").Value = "No"
$ws.Range("F7").Value = "No"
$ws.Range("F8").Value = "No"
$ws.Range("F9").Value = "No"
$ws.Range("F10").Value = "No"
$ws.Range("F11").Value = "No"
$ws.Range("F12").Value = "No"
$ws.Range("F13").Value = "No"
$ws.Range("F14").Value = "No"
$ws.Range("F15").Value = "No"
$ws.Range("F17").Value = "No"
$ws.Range("F18").Value = "No"
$ws.Range("F19").Value = "No"
$ws.Range("F20").Value = "No"
$ws.Range("F21").Value = "No"

# Update the view: scroll back to top-left and move the active selection to F22.
$ws.Range("A1").Select() | Out-Null
$ws.Range("F22").Select() | Out-Null
